$d = $word.ActiveDocument

$xmlTail = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>tail –f nexus.log</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$xmlCountdown = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>(10..1).each{ log.info '--&gt; ' + it }</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$xmlBlastoff = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>log.info 'Blast off!'</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$xmlLogVsPrintln = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>log vs println</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$xmlPrintlnOnly = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>p</w:t></w:r><w:r><w:t>rintln only goes to the task log. This output does not appear in the nexus log in the browser.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$xmlDefList1 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>def list = ['monday', 'wednesday', 'chocolate', 'friday']</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$xmlPrintlnBlank = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>println ______________________________________</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$xmlDefList2 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>def list = ['monday', 'wednesday', 'chocolate', 'friday']</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$xmlGrepPlusNewSection = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>println list.grep { it.endsWith 'day' }</w:t></w:r></w:p><w:p/><w:p/><w:p><w:r><w:br w:type="page"/></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Creating a user</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Task: </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Add four more lines of code so Owen, Sam, </w:t></w:r><w:r><w:t>Sophia</w:t></w:r><w:r><w:t xml:space="preserve"> and Daisy get accounts. </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Solution:</w:t></w:r></w:p><w:p><w:r><w:t>import jenkins.model.Jenkins</w:t></w:r></w:p><w:p><w:r><w:t>def instance = Jenkins.getInstance()</w:t></w:r></w:p><w:p><w:r><w:t>def realm = Jenkins.getInstance().securityRealm</w:t></w:r></w:p><w:p><w:r><w:t>realm.createAccount('olivia', 'olivia')</w:t></w:r></w:p><w:p><w:r><w:t>realm.createAccount('owen', 'owen')</w:t></w:r></w:p><w:p><w:r><w:t>realm.createAccount('sam', 'sam')</w:t></w:r></w:p><w:p><w:r><w:t>realm.createAccount(</w:t></w:r><w:r><w:t>'sophia</w:t></w:r><w:r><w:t>', 's</w:t></w:r><w:r><w:t>ophia</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>')</w:t></w:r></w:p><w:p><w:r><w:t>realm.createAccount('daisy', 'daisy')</w:t></w:r></w:p><w:p><w:r><w:t>instance.save()</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

function Replace-ParaXml($paraIndex, $xml) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.InsertXML($xml)
}

Replace-ParaXml 9  $xmlTail
Replace-ParaXml 50 $xmlCountdown
Replace-ParaXml 51 $xmlBlastoff
Replace-ParaXml 54 $xmlLogVsPrintln
Replace-ParaXml 59 $xmlPrintlnOnly
Replace-ParaXml 65 $xmlDefList1
Replace-ParaXml 66 $xmlPrintlnBlank
Replace-ParaXml 69 $xmlDefList2
Replace-ParaXml 70 $xmlGrepPlusNewSection

Write-Host "done"
